# ExcelTestSource plugin: allow defining first step in test case row,
# improve flexibility.
#
# The sample data modelled one "first test step" per test case as its own
# row (row 11, 12, ...) below the test-case row. This edit folds the first
# step's Action/Expected (and its step number) directly into the test-case
# row itself (row 10), and removes the now-redundant extra row, so a test
# case's first step can be defined inline on the test case row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pull the "step 1" values (D11:F11 = step number / action / expected)
# up into the test case row (row 10), which previously only had the
# test-case columns (A, C, G, I, J, L, M, N) populated.
$ws.Range("D10").Value = $ws.Range("D11").Value2
$ws.Range("E10").Value = $ws.Range("E11").Value2
$ws.Range("F10").Value = $ws.Range("F11").Value2

# The old "step 1" row (11) is now redundant - delete it, shifting the old
# "step 2" row (12) up to become row 11.
$ws.Rows("11:11").Delete()

# The table / autofilter range, used range dimension, etc. shrink by one
# row automatically as a result of the delete above.

# Column C ("Title") and D ("Test Step") grew to fit the values that moved
# into row 10 - best-fit the widths to the new content.
$ws.Columns.Item(3).ColumnWidth = 28.5
$ws.Columns.Item(4).ColumnWidth = 10.666666666666666

# Update the active selection to reflect where the edit was made.
$ws.Range("A10").Select()
